# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
#
# NOTE: Excel's Range/Cells ".Value" setter auto-detects dates
# ("2026-01-28") and percentages ("91.0%") and silently converts them to
# numeric/date values instead of keeping the literal text. To preserve the
# source data exactly as plain text (matching the rest of the log), the
# target cells are switched to the "@" (Text) number format immediately
# before the writes, then restored to "General" afterwards.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: append rows 188-200
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    @("2026-01-28","15:03:12","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:14","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:19","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:24","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:29","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:34","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:39","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:44","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:49","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:54","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:03:59","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:04:04","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:04:09","15:00","Bathroom","No Motion","Inactive")
)
$startRow = 188
$endRow = $startRow + $rows.Count - 1
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
$ws.Range("A$startRow`:A$endRow").NumberFormat = "General"

# ---------------------------------------------------------------------
# Humidity sheet: append rows 183-193
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    @("2026-01-28","15:03:11","15:00","Bathroom","91.0%","Active"),
    @("2026-01-28","15:03:13","15:00","Bathroom","90.2%","Active"),
    @("2026-01-28","15:03:15","15:00","Bathroom","90.4%","Active"),
    @("2026-01-28","15:03:19","15:00","Bathroom","89.5%","Active"),
    @("2026-01-28","15:03:23","15:00","Bathroom","90.1%","Active"),
    @("2026-01-28","15:03:31","15:00","Bathroom","89.0%","Active"),
    @("2026-01-28","15:03:35","15:00","Bathroom","89.8%","Active"),
    @("2026-01-28","15:03:39","15:00","Bathroom","88.6%","Active"),
    @("2026-01-28","15:03:47","15:00","Bathroom","89.4%","Active"),
    @("2026-01-28","15:03:55","15:00","Bathroom","89.2%","Active"),
    @("2026-01-28","15:04:07","15:00","Bathroom","89.1%","Active")
)
$startRow = 183
$endRow = $startRow + $rows.Count - 1
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("E$startRow`:E$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
$ws.Range("A$startRow`:A$endRow").NumberFormat = "General"
$ws.Range("E$startRow`:E$endRow").NumberFormat = "General"

# ---------------------------------------------------------------------
# Temperature sheet: append rows 183-193
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
    @("2026-01-28","15:03:12","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:13","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:16","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:20","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:24","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:32","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:36","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:40","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:48","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:03:56","15:00","Bathroom","23.0C","Active"),
    @("2026-01-28","15:04:08","15:00","Bathroom","23.0C","Active")
)
$startRow = 183
$endRow = $startRow + $rows.Count - 1
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
$ws.Range("A$startRow`:A$endRow").NumberFormat = "General"
